# Auto-generated edit script: applies numeric corrections to the Profits tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 299.68182
$ws.Range("I92").Value = 250.5
$ws.Range("K92").Value = 250.5
$ws.Range("M92").Value = 997.5
$ws.Range("H101").Value = 2540
$ws.Range("I101").Value = 300
$ws.Range("J101").Value = 3100
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 9300
$ws.Range("M101").Value = 722
$ws.Range("N101").Value = -12544
$ws.Range("H103").Value = 1203.0834
$ws.Range("I103").Value = 1029.75
$ws.Range("J103").Value = 1549.75
$ws.Range("K103").Value = 3089.25
$ws.Range("L103").Value = 4649.25
$ws.Range("M103").Value = -2503.25
$ws.Range("N103").Value = -5821.25
$ws.Range("H106").Value = 2933.3333
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 690.7857
$ws.Range("I107").Value = 690.7857
$ws.Range("K107").Value = 690.7857
$ws.Range("M107").Value = 1229.2143
$ws.Range("H108").Value = 40684
$ws.Range("J108").Value = 40684
$ws.Range("L108").Value = 40684
$ws.Range("N108").Value = -48364
$ws.Range("H129").Value = 1284.3793
$ws.Range("I129").Value = 593
$ws.Range("J129").Value = 1395
$ws.Range("K129").Value = 1779
$ws.Range("L129").Value = 4185
$ws.Range("M129").Value = 3221
$ws.Range("N129").Value = -14185
$ws.Range("H132").Value = 348753.2
$ws.Range("I132").Value = 348753.2
$ws.Range("K132").Value = 1046259.6
$ws.Range("M132").Value = -1043729.6
$ws.Range("H137").Value = 4417.244
$ws.Range("I137").Value = 5402.36
$ws.Range("J137").Value = 2878
$ws.Range("K137").Value = 16207.08
$ws.Range("L137").Value = 8634
$ws.Range("M137").Value = -13657.08
$ws.Range("N137").Value = -13734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4902460.5
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 5882833
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 5882833
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -5883059
$ws.Range("H32").Value = 10538.867
$ws.Range("I32").Value = 4460.891
$ws.Range("K32").Value = 4460.891
$ws.Range("M32").Value = -4173.891
$ws.Range("H45").Value = 38910.258
$ws.Range("I45").Value = 57097.945
$ws.Range("J45").Value = 2534.889
$ws.Range("K45").Value = 57097.945
$ws.Range("L45").Value = 2534.889
$ws.Range("M45").Value = -56720.945
$ws.Range("N45").Value = -3288.889
$ws.Range("H97").Value = 618.5263
$ws.Range("J97").Value = 650.2222
$ws.Range("L97").Value = 650.2222
$ws.Range("N97").Value = -1642.2222
$ws.Range("H102").Value = 2308.75
$ws.Range("I102").Value = 2210
$ws.Range("K102").Value = 2210
$ws.Range("M102").Value = -588
$ws.Range("H110").Value = 42922
$ws.Range("I110").Value = 876.25
$ws.Range("K110").Value = 876.25
$ws.Range("M110").Value = 1168.75
$ws.Range("H116").Value = 4902460.5
$ws.Range("I116").Value = 600
$ws.Range("J116").Value = 5882833
$ws.Range("K116").Value = 600
$ws.Range("L116").Value = 5882833
$ws.Range("M116").Value = 1694
$ws.Range("N116").Value = -5887421
$ws.Range("H132").Value = 1151167.8
$ws.Range("I132").Value = 1438683.9
$ws.Range("J132").Value = 1103.1666
$ws.Range("K132").Value = 4316051.699999999
$ws.Range("L132").Value = 3309.4998
$ws.Range("M132").Value = -4313521.699999999
$ws.Range("N132").Value = -8369.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4902460.5
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 5882833
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 5882833
$ws.Range("M3").Value = -486
$ws.Range("N3").Value = -5883061
$ws.Range("H94").Value = 1172.3478
$ws.Range("I94").Value = 1062.1177
$ws.Range("J94").Value = 1484.6666
$ws.Range("K94").Value = 1062.1177
$ws.Range("L94").Value = 1484.6666
$ws.Range("M94").Value = -611.1177
$ws.Range("N94").Value = -2386.6666
$ws.Range("H99").Value = 1230
$ws.Range("I99").Value = 1341.6666
$ws.Range("J99").Value = 895
$ws.Range("K99").Value = 1341.6666
$ws.Range("L99").Value = 895
$ws.Range("M99").Value = 156.3334
$ws.Range("N99").Value = -3891
$ws.Range("H105").Value = 2037.4375
$ws.Range("I105").Value = 2126.2727
$ws.Range("J105").Value = 1842
$ws.Range("K105").Value = 2126.2727
$ws.Range("L105").Value = 1842
$ws.Range("M105").Value = -379.2727
$ws.Range("N105").Value = -5336
$ws.Range("H107").Value = 672660
$ws.Range("I107").Value = 830015.3
$ws.Range("J107").Value = 3900
$ws.Range("K107").Value = 830015.3
$ws.Range("L107").Value = 3900
$ws.Range("M107").Value = -828095.3
$ws.Range("N107").Value = -7740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3317.4324
$ws.Range("I58").Value = 2061.4
$ws.Range("J58").Value = 4173.8184
$ws.Range("K58").Value = 2061.4
$ws.Range("L58").Value = 4173.8184
$ws.Range("M58").Value = -1858.4
$ws.Range("N58").Value = -4579.8184
$ws.Range("H94").Value = 1217.7693
$ws.Range("I94").Value = 1103
$ws.Range("J94").Value = 1268.7778
$ws.Range("K94").Value = 1103
$ws.Range("L94").Value = 1268.7778
$ws.Range("M94").Value = -652
$ws.Range("N94").Value = -2170.7778
$ws.Range("H105").Value = 920
$ws.Range("I105").Value = 920
$ws.Range("K105").Value = 920
$ws.Range("M105").Value = 827
$ws.Range("H107").Value = 1872.2333
$ws.Range("I107").Value = 1812.0588
$ws.Range("J107").Value = 1950.9231
$ws.Range("K107").Value = 1812.0588
$ws.Range("L107").Value = 1950.9231
$ws.Range("M107").Value = 107.9412
$ws.Range("N107").Value = -5790.9231
$ws.Range("H134").Value = 2055.3333
$ws.Range("I134").Value = 2131.147
$ws.Range("J134").Value = 1821
$ws.Range("K134").Value = 6393.441
$ws.Range("L134").Value = 5463
$ws.Range("M134").Value = -3858.441
$ws.Range("N134").Value = -10533
$ws.Range("H136").Value = 3317.4324
$ws.Range("I136").Value = 2061.4
$ws.Range("J136").Value = 4173.8184
$ws.Range("K136").Value = 6184.200000000001
$ws.Range("L136").Value = 12521.4552
$ws.Range("M136").Value = -3634.200000000001
$ws.Range("N136").Value = -17621.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 456.30356
$ws.Range("I5").Value = 459.12
$ws.Range("J5").Value = 454.03226
$ws.Range("K5").Value = 1377.36
$ws.Range("L5").Value = 1362.09678
$ws.Range("M5").Value = -1265.36
$ws.Range("N5").Value = -1586.09678
$ws.Range("H6").Value = 61.8
$ws.Range("I6").Value = 61.8
$ws.Range("K6").Value = 185.4
$ws.Range("M6").Value = -72.39999999999998
$ws.Range("H68").Value = 957.4267
$ws.Range("I68").Value = 754.6667
$ws.Range("J68").Value = 1144.5897
$ws.Range("K68").Value = 2264.0001
$ws.Range("L68").Value = 3433.7691
$ws.Range("M68").Value = -1453.0001
$ws.Range("N68").Value = -5055.7691
$ws.Range("H71").Value = 957.4267
$ws.Range("I71").Value = 754.6667
$ws.Range("J71").Value = 1144.5897
$ws.Range("K71").Value = 6792.0003
$ws.Range("L71").Value = 10301.3073
$ws.Range("M71").Value = -2736.0003
$ws.Range("N71").Value = -18413.3073
$ws.Range("H98").Value = 1500
$ws.Range("I98").Value = 1500
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4500
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = -3002
$ws.Range("H107").Value = 49573.39
$ws.Range("I107").Value = 31894.062
$ws.Range("J107").Value = 112433.22
$ws.Range("K107").Value = 95682.186
$ws.Range("L107").Value = 337299.66
$ws.Range("M107").Value = -93762.186
$ws.Range("N107").Value = -341139.66
$ws.Range("H131").Value = 1668601.6
$ws.Range("I131").Value = 930.5294
$ws.Range("J131").Value = 2327913.2
$ws.Range("K131").Value = 2791.5882
$ws.Range("L131").Value = 6983739.600000001
$ws.Range("M131").Value = 2248.4118
$ws.Range("N131").Value = -6993819.600000001
$ws.Range("H135").Value = 456.30356
$ws.Range("I135").Value = 459.12
$ws.Range("J135").Value = 454.03226
$ws.Range("K135").Value = 4132.08
$ws.Range("L135").Value = 4086.29034
$ws.Range("M135").Value = -1597.08
$ws.Range("N135").Value = -9156.29034

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 405
$ws.Range("I13").Value = 405
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 405
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = -266
$ws.Range("H97").Value = 1233.72
$ws.Range("I97").Value = 1231.8695
$ws.Range("J97").Value = 1255
$ws.Range("K97").Value = 1231.8695
$ws.Range("L97").Value = 1255
$ws.Range("M97").Value = -735.8695
$ws.Range("N97").Value = -2247

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1939.7878
$ws.Range("J93").Value = 2468.8667
$ws.Range("L93").Value = 2468.8667
$ws.Range("N93").Value = -4964.8667
$ws.Range("H100").Value = 2200.4443
$ws.Range("I100").Value = 2150
$ws.Range("J100").Value = 2301.3333
$ws.Range("K100").Value = 2150
$ws.Range("L100").Value = 2301.3333
$ws.Range("M100").Value = -1609
$ws.Range("N100").Value = -3383.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 14001618
$ws.Range("I17").Value = 23334000
$ws.Range("J17").Value = 3045
$ws.Range("K17").Value = 23334000
$ws.Range("L17").Value = 3045
$ws.Range("M17").Value = -23333828
$ws.Range("N17").Value = -3389
$ws.Range("H96").Value = 14287186
$ws.Range("I96").Value = 25001124
$ws.Range("J96").Value = 1933.3334
$ws.Range("K96").Value = 25001124
$ws.Range("L96").Value = 1933.3334
$ws.Range("M96").Value = -24999751
$ws.Range("N96").Value = -4679.3334
$ws.Range("H100").Value = 964.05884
$ws.Range("I100").Value = 1435.75
$ws.Range("J100").Value = 544.7778
$ws.Range("K100").Value = 2871.5
$ws.Range("L100").Value = 1089.5556
$ws.Range("M100").Value = -2330.5
$ws.Range("N100").Value = -2171.5556
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H113").Value = 535.3684
$ws.Range("I113").Value = 416.54544
$ws.Range("K113").Value = 1249.63632
$ws.Range("M113").Value = 920.3636799999999
$ws.Range("H132").Value = 5439
$ws.Range("I132").Value = 5786.59
$ws.Range("J132").Value = 2050
$ws.Range("K132").Value = 17359.77
$ws.Range("L132").Value = 6150
$ws.Range("M132").Value = -14829.77
$ws.Range("N132").Value = -11210
